$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 4 (shifts old row4 ING data down to row5) ---
$ws.Rows("4").Insert()

# --- Row 2 updates ---
$ws.Range("B2").Value = "'3"
$ws.Range("D2").Value = -0.0204
$ws.Range("E2").Value = -0.224
$ws.Range("F2").Value = -0.09914
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2347.4
$ws.Range("L2").Value = 0.09983158682634731
$ws.Range("M2").Value = 762.41
$ws.Range("N2").Value = 0.01631671649648908
$ws.Range("O2").Value = 0.3247891283973758
$ws.Range("P2").Value = 757.73
$ws.Range("Q2").Value = 0.01621655748335499
$ws.Range("R2").Value = 0.3227954332452926
$ws.Range("S2").Value = 4.679999999999999
$ws.Range("T2").Value = 0.006138429453968336
$ws.Range("U2").Value = 202176.9
$ws.Range("V2").Value = 4.326888628741784
$ws.Range("W2").Value = 0.01729750430841897
$ws.Range("X2").Value = 0.142979431556923
$ws.Range("Y2").Value = -0.125681927248504
$ws.Range("Z2").Value = 0.0669975475984444
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.03274378739066153
$ws.Range("AC2").Value = -0.03274378739066153
$ws.Range("AD2").Value = 411900.9
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 411900.9
$ws.Range("AG2").Value = 209724
$ws.Range("AH2").Value = 0.8981182077097142
$ws.Range("AI2").Value = 0.8193644498852017
$ws.Range("AJ2").Value = 0.8177977981647083
$ws.Range("AK2").Value = 0.6978450468820344
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0

# --- Row 3 updates ---
$ws.Range("B3").Value = 'Van Lanschot Kempen N.V. (ENXTAM:VLK)'
$ws.Range("D3").Value = -0.0052
$ws.Range("E3").Value = -0.224
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 27.1
$ws.Range("L3").Value = 0.05126749905410519
$ws.Range("M3").Value = 12.29
$ws.Range("N3").Value = 0.01170699180796342
$ws.Range("O3").Value = 0.4535055350553505
$ws.Range("P3").Value = 7.61
$ws.Range("Q3").Value = 0.007248999809487522
$ws.Range("R3").Value = 0.2808118081180812
$ws.Range("S3").Value = 4.679999999999999
$ws.Range("T3").Value = 0.3807973962571196
$ws.Range("U3").Value = 2371.8
$ws.Range("V3").Value = 2.259287483330159
$ws.Range("W3").Value = 0.01729750430841897
$ws.Range("X3").Value = 0.09610557732932715
$ws.Range("Y3").Value = -0.07880807302090817
$ws.Range("Z3").Value = 0.1547605105984307
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.03148349401109913
$ws.Range("AC3").Value = -0.03148349401109913
$ws.Range("AD3").Value = 3971.1
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 3971.1
$ws.Range("AG3").Value = 1599.3
$ws.Range("AH3").Value = 0.7909139795654166
$ws.Range("AI3").Value = 0.7268018595116951
$ws.Range("AJ3").Value = 0.6037144690649655
$ws.Range("AK3").Value = 0.5172380336351875
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0
# remove F3 (no longer present)
$ws.Range("F3").ClearContents()

# --- Row 4 (new ABN AMRO) ---
$ws.Range("A4").Value = 'Netherlands'
$ws.Range("B4").Value = 'ABN AMRO Bank N.V. (ENXTAM:ABN)'
$ws.Range("C4").Value = 'Bank (Money Center)'
$ws.Range("D4").Value = -0.0614
$ws.Range("E4").Value = -0.362
$ws.Range("F4").Value = -0.198
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 254.8
$ws.Range("L4").Value = 0.03728307629276287
$ws.Range("M4").Value = 750.12
$ws.Range("N4").Value = 0.08134203734628816
$ws.Range("O4").Value = 2.943956043956044
$ws.Range("P4").Value = 750.12
$ws.Range("Q4").Value = 0.08134203734628816
$ws.Range("R4").Value = 2.943956043956044
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 69367.1
$ws.Range("V4").Value = 7.522078119239195
$ws.Range("W4").Value = 0.01113013698630137
$ws.Range("X4").Value = 0.3434678308875967
$ws.Range("Y4").Value = -0.3323376939012953
$ws.Range("Z4").Value = 0.05422965435756721
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.03274378739066153
$ws.Range("AC4").Value = -0.03274378739066153
$ws.Range("AD4").Value = 169326
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 169326
$ws.Range("AG4").Value = 99958.9
$ws.Range("AH4").Value = 0.9483510858156752
$ws.Range("AI4").Value = 0.8737737825550733
$ws.Range("AJ4").Value = 0.9155363539526674
$ws.Range("AK4").Value = 0.8033996169423059
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0

# --- Row 5 (ING updated values) ---
$ws.Range("D5").Value = -0.0204
$ws.Range("E5").Value = -0.166
$ws.Range("F5").Value = -0.00028
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2065.5
$ws.Range("L5").Value = 0.1278884018129133
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("U5").Value = 130438
$ws.Range("V5").Value = 3.578143473573618
$ws.Range("W5").Value = 0.03539900358360283
$ws.Range("X5").Value = 0.142979431556923
$ws.Range("Y5").Value = -0.1075804279733201
$ws.Range("Z5").Value = 0.07290793921359028
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.03367471697282425
$ws.Range("AC5").Value = -0.03367471697282425
$ws.Range("AD5").Value = 238603.8
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 238603.8
$ws.Range("AG5").Value = 108165.8
$ws.Range("AH5").Value = 0.8674675404705701
$ws.Range("AI5").Value = 0.7862853715682947
$ws.Range("AJ5").Value = 0.7479316470278294
$ws.Range("AK5").Value = 0.6251671781711834
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0
# remove T5, AN5, AP5 leftover from shifted row (need to ensure cleared)
$ws.Range("T5").ClearContents()
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()

# --- clear AN2,AO2,AP2,AQ2 leftover (already absent in before but ensure) and AN4 etc from row3 ---
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
